$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.161.45"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "1.680.63"
$ws.Range("E3").Value = "  +0.43%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.25"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("E6").Value = "  +0.24%  "

$ws.Range("E8").Value = "  +2.01%  "

$ws.Range("E9").Value = "  +5.61%  "

$ws.Range("E10").Value = "  +0.70%  "

$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").Value = "1.916.98"
$ws.Range("E12").Value = "  +0.40%  "

$ws.Range("D13").Value = "1.690.35"
$ws.Range("E13").Value = "  +0.98%  "

$ws.Range("E14").Value = "  +1.60%  "

$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.24"
$ws.Range("E16").Value = "  +0.98%  "

$ws.Range("D17").Value = "27.154.43"
$ws.Range("E17").Value = "  +0.79%  "

$ws.Range("E18").Value = "  +1.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.08"
$ws.Range("E19").Value = "  -0.08%  "

$ws.Range("D20").Value = "0.0₃0745"
$ws.Range("E20").Value = "  +1.47%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.52"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.46"
$ws.Range("E23").Value = "  +3.02%  "

$ws.Range("E24").Value = "  -2.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.91"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("E27").Value = "  +2.19%  "

$ws.Range("E28").Value = "  +0.61%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("E30").Value = "  +0.40%  "

$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("D32").Value = "1.560.55"
$ws.Range("E32").Value = "  +5.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.36"
$ws.Range("E33").Value = "  +1.27%  "

$ws.Range("E34").Value = "  +2.71%  "

$ws.Range("E35").Value = "  +0.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.603"
$ws.Range("E36").Value = "  +3.30%  "

$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.935"
$ws.Range("E38").Value = "  +4.53%  "

$ws.Range("E39").Value = "  +2.73%  "

$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.32"
$ws.Range("E41").Value = "  +2.86%  "

$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("E43").Value = "  -3.05%  "

$ws.Range("E44").Value = "  -2.39%  "

$ws.Range("D45").Value = "1.825.68"
$ws.Range("E45").Value = "  +0.64%  "

$ws.Range("E46").Value = "  +0.61%  "

$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("E48").Value = "  +3.96%  "

$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  +0.48%  "

$ws.Range("E50").Value = "  +2.62%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.17"
$ws.Range("E51").Value = "  +5.71%  "
